$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header-style labels (written as text -> become shared strings)
$ws.Range("A1").Value = "One"
$ws.Range("B1").Value = "Two"
$ws.Range("C1").Value = "great"
$ws.Range("D1").Value = "test"
$ws.Range("E1").Value = "Thirty"

# Row 2
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 45
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = 30

# Row 3
$ws.Range("A3").Value = 45
$ws.Range("B3").Value = 54
$ws.Range("C3").Value = 787
$ws.Range("D3").Value = 67
$ws.Range("E3").Value = 343

# Page orientation -> portrait
$ws.PageSetup.Orientation = 1

# Leave the selection on E1, matching the saved view state
$ws.Range("E1").Select()
